# Update countries & provincias Spain
# - Swap the "El Salvador" / "Burkina Faso" rows (109/110) and update their daily figures
# - Swap the "Brunei" / "Guayana Francesa" rows (152/153) and update their daily figures
# - Refresh the "Datos actualizados" timestamp string
# - Refresh the daily COVID figures for several other countries (Estados Unidos,
#   Alemania, Brasil, Canada, India, Mali, Guadalupe)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country labels: El Salvador (row 109) <-> Burkina Faso (row 110) ---
$ws.Range("A109").Value = "Burkina Faso"
$ws.Range("A110").Value = "El Salvador"

# --- Swap country labels: Brunei (row 152) <-> Guayana Francesa (row 153) ---
$ws.Range("A152").Value = "Guayana Francesa"
$ws.Range("A153").Value = "Brunei"

# --- Timestamp string update ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 00:04"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1318686
$ws.Range("C4").Value = 26063
$ws.Range("D4").Value = 222008
$ws.Range("E4").Value = 1018180
$ws.Range("G4").Value = 1570
$ws.Range("H4").Value = 78498

# --- Row 10: Alemania ---
$ws.Range("B10").Value = 170678
$ws.Range("C10").Value = 1248
$ws.Range("E10").Value = 21468
$ws.Range("G10").Value = 118
$ws.Range("H10").Value = 7510

# --- Row 11: Brasil ---
$ws.Range("B11").Value = 145328
$ws.Range("C11").Value = 9635
$ws.Range("E11").Value = 80081
$ws.Range("G11").Value = 709
$ws.Range("H11").Value = 9897

# --- Row 15: Canada ---
$ws.Range("B15").Value = 66326
$ws.Range("C15").Value = 1404
$ws.Range("D15").Value = 29948
$ws.Range("E15").Value = 31811

# --- Row 17: India ---
$ws.Range("D17").Value = 17887
$ws.Range("E17").Value = 39821

# --- Row 109: now Burkina Faso (new figures) ---
$ws.Range("B109").Value = 744
$ws.Range("C109").Value = 8
$ws.Range("D109").Value = 566
$ws.Range("E109").Value = 130
$ws.Range("F109").Value = 0
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 48

# --- Row 110: now El Salvador (figures that used to belong to the old row 109) ---
$ws.Range("B110").Value = 742
$ws.Range("C110").Value = 47
$ws.Range("D110").Value = 257
$ws.Range("E110").Value = 469
$ws.Range("F110").Value = 4
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = 16

# --- Row 113: Mali ---
$ws.Range("B113").Value = 668
$ws.Range("C113").Value = 18
$ws.Range("D113").Value = 285
$ws.Range("E113").Value = 348
$ws.Range("G113").Value = 3
$ws.Range("H113").Value = 35

# --- Row 149: Guadalupe ---
$ws.Range("B149").Value = 154
$ws.Range("C149").Value = 1
$ws.Range("E149").Value = 37

# --- Row 152: now Guayana Francesa (new figures) ---
$ws.Range("C152").Value = 3
$ws.Range("D152").Value = 113
$ws.Range("E152").Value = 27
$ws.Range("F152").Value = 0

# --- Row 153: now Brunei (figures that used to belong to the old row 152) ---
$ws.Range("B153").Value = 141
$ws.Range("D153").Value = 132
$ws.Range("E153").Value = 8
$ws.Range("F153").Value = 2
